$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextCell "D2" '67.637.80'
Set-TextCell "E2" '  -0.17%  '

# Row 3
Set-TextCell "D3" '3.769.88'
Set-TextCell "E3" '  -0.80%  '

# Row 4
Set-TextCell "E4" '  +0.07%  '

# Row 5
Set-TextCell "D5" '598.21'
Set-TextCell "E5" '  +0.24%  '

# Row 6
Set-TextCell "D6" '162.98'
Set-TextCell "E6" '  -2.55%  '

# Row 7
Set-TextCell "D7" '3.767.54'
Set-TextCell "E7" '  -0.91%  '

# Row 8
Set-TextCell "E8" '  -0.11%  '

# Row 9
Set-TextCell "D9" '0.512'
Set-TextCell "E9" '  -1.26%  '

# Row 10
Set-TextCell "E10" '  -3.08%  '

# Row 11
Set-TextCell "E11" '  -1.14%  '

# Row 12
Set-TextCell "E12" '  +3.94%  '

# Row 13
Set-TextCell "D13" '0.0000244'
Set-TextCell "E13" '  -3.87%  '

# Row 14
Set-TextCell "D14" '35.16'
Set-TextCell "E14" '  -2.10%  '

# Row 15
Set-TextCell "D15" '4.400.01'
Set-TextCell "E15" '  -0.96%  '

# Row 16
Set-TextCell "D16" '3.776.99'
Set-TextCell "E16" '  -0.56%  '

# Row 17
Set-TextCell "D17" '67.698.80'
Set-TextCell "E17" '  -0.15%  '

# Row 18
Set-TextCell "E18" '  -1.55%  '

# Row 19
Set-TextCell "D19" '0.115'
Set-TextCell "E19" '  +1.73%  '

# Row 20
Set-TextCell "E20" '  -1.35%  '

# Row 21
Set-TextCell "D21" '456.91'
Set-TextCell "E21" '  -0.99%  '

# Row 22
Set-TextCell "D22" '9.47'
Set-TextCell "E22" '  -4.48%  '

# Row 23
Set-TextCell "D23" '0.692'
Set-TextCell "E23" '  -1.09%  '

# Row 24
Set-TextCell "D24" '82.67'
Set-TextCell "E24" '  -0.73%  '

# Row 25
Set-TextCell "E25" '  -6.06%  '

# Row 26
Set-TextCell "E26" '  -1.94%  '

# Row 27
Set-TextCell "B27" 'Fetch.AI'
Set-TextCell "C27" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell "D27" '2.08'
Set-TextCell "E27" '  -1.19%  '

# Row 28
Set-TextCell "B28" 'Dai'
Set-TextCell "C28" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell "D28" '1.00'
Set-TextCell "E28" '  +0.03%  '

# Row 29
Set-TextCell "D29" '9.83'
Set-TextCell "E29" '  -1.67%  '

# Row 30
Set-TextCell "D30" '3.918.13'
Set-TextCell "E30" '  -0.71%  '

# Row 31
Set-TextCell "E31" '  -0.91%  '

# Row 32
Set-TextCell "D32" '7.23'
Set-TextCell "E32" '  -1.65%  '

# Row 33
Set-TextCell "D33" '2.58'
Set-TextCell "E33" '  -6.93%  '

# Row 34
Set-TextCell "D34" '28.80'
Set-TextCell "E34" '  -2.54%  '

# Row 35
Set-TextCell "E35" '  +0.06%  '

# Row 36
Set-TextCell "D36" '8.93'
Set-TextCell "E36" '  -1.33%  '

# Row 37
Set-TextCell "D37" '0.0987'
Set-TextCell "E37" '  -1.42%  '

# Row 38
Set-TextCell "E38" '  +2.52%  '

# Row 39
Set-TextCell "D39" '5.76'
Set-TextCell "E39" '  -0.27%  '

# Row 40
Set-TextCell "D40" '0.976'
Set-TextCell "E40" '  -2.26%  '

# Row 41
Set-TextCell "B41" 'FirstDigitalUSD'
Set-TextCell "C41" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell "D41" '1.00'
Set-TextCell "E41" '  +0.04%  '

# Row 42
Set-TextCell "B42" 'dogwifhat'
Set-TextCell "C42" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell "D42" '3.14'
Set-TextCell "E42" '  -6.03%  '

# Row 43
Set-TextCell "E43" '  +0.03%  '

# Row 44
Set-TextCell "D44" '47.19'
Set-TextCell "E44" '  -1.88%  '

# Row 45
Set-TextCell "D45" '43.03'
Set-TextCell "E45" '  +0.49%  '

# Row 46
Set-TextCell "D46" '152.11'
Set-TextCell "E46" '  +3.18%  '

# Row 47
Set-TextCell "D47" '0.294'
Set-TextCell "E47" '  -2.56%  '

# Row 48
Set-TextCell "D48" '1.36'
Set-TextCell "E48" '  +0.27%  '

# Row 49
Set-TextCell "D49" '8.26'
Set-TextCell "E49" '  -1.04%  '

# Row 50
Set-TextCell "B50" 'Stacks'
Set-TextCell "C50" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell "D50" '1.84'
Set-TextCell "E50" '  -0.35%  '

# Row 51
Set-TextCell "B51" 'Bittensor'
Set-TextCell "C51" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell "D51" '385.29'
Set-TextCell "E51" '  -2.64%  '
